$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (already sorted by expected_profit_gold, column E, descending)
# columns: name, id, highestBuyPrice, lowestSellPrice, expected_profit_gold, roi_pct,
#          status, T_days, avg_units_day, bursty, zero_streak7, prob_exit_3d, reason
$data = @(
    ,("Shadow Chevalier Torch Skin", 27416, 868790, 4000000, 253.121, 291.3488875332359, "DEADLIKE", "inf", 0, $false, 6, 0, "6 consecutive zero-unit days in last 7")
    ,("Winter''s Heart Infusion", 24511, 1200349, 1985002, 48.69027, 40.5634278030806, "HOLD_RISK", 3970004, 0, $false, 1, 0.000000755666459406612, "")
    ,("Crustacea", 19937, 1350000, 1890000, 25.65, 19, "OK_SLOW", 5, 0.235, $false, 2, 0.505852993582407, "")
    ,("Polysaturating Reverberating Infusion (Gray)", 25627, 1470530, 1989893, 22.08790500000001, 15.02037020666019, "DEADLIKE", "inf", 0, $false, 4, 0, "4 consecutive zero-unit days in last 7")
    ,("Phospholuminescent Infusion", 24480, 1390711, 1849999, 18.17881499999999, 13.07159790927086, "DEADLIKE", "inf", 0, $false, 3, 0, "3 consecutive zero-unit days in last 7")
    ,("Peerless Infusion", 25865, 1027081, 1349992, 12.04122, 11.72372967662725, "OK_LIQUID", 2, 0.544, $false, 2, 0.804581071730771, "")
    ,("Winter''s Heart Infusion", 24517, 955004, 1199999, 6.499515000000002, 6.805746363366019, "OK_LIQUID", 2, 0.738, $false, 1, 0.8905749669380867, "")
    ,("Ghostly Infusion", 24206, 523241, 679907, 5.467994999999996, 10.45024185795837, "OK_LIQUID", 3, 0.497, $false, 2, 0.7750489414166475, "")
    ,("Ghostly Infusion", 24216, 518379, 669889, 5.102665000000003, 9.843502533860365, "HOLD_RISK", 365, 0.003, $false, 2, 0.008204093178705296, "")
    ,("Ghostly Infusion", 24207, 640010, 779791, 2.281234999999997, 3.564373994156337, "NO_FLOW", "inf", 0, $false, 2, 0, "Trading observed, but effective daily sell-through rate collapsed to ~0 due to sparse or volatile recent activity")
    ,("Zojja''s Berserker Insignia", 20887, 205601, 260000, 1.5399, 7.489749563474886, "OK_LIQUID", 1, 14.58, $false, 2, 1, "")
    ,("Zojja''s Berserker Inscription", 20865, 250543, 309543, 1.256854999999999, 5.01652410963387, "DEADLIKE", "inf", 0, $false, 3, 0, "3 consecutive zero-unit days in last 7")
    ,("Vintage Black Lion Weapon Box", 26458, 420001, 508725, 1.241525, 2.956004866655079, "DEADLIKE", "inf", 0, $false, 3, 0, "3 consecutive zero-unit days in last 7")
    ,("Jade Bot Core: Tier 10", 26651, 447630, 472000, -4.643, -10.37240578156066, "OK_LIQUID", 1, 31.375, $false, 1, 1, "")
    ,("Mystic Facet", 27392, 725017, 794999, -4.926784999999998, -6.79540617668275, "DEADLIKE", "inf", 0, $false, 4, 0, "4 consecutive zero-unit days in last 7")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws.Cells.Item($rowIndex, 9).Value = $row[8]
    $ws.Cells.Item($rowIndex, 10).Value = $row[9]
    $ws.Cells.Item($rowIndex, 11).Value = $row[10]
    $ws.Cells.Item($rowIndex, 12).Value = $row[11]
    $ws.Cells.Item($rowIndex, 13).Value = $row[12]
    $rowIndex = $rowIndex + 1
}
